$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before D to make room for the newest two quarters
$ws.Columns("D:E").Insert()

# Copy number formatting from column F (the old column D, now shifted) into new D:E
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new quarters of data, and correct a few restated prior-quarter figures
$ws.Range("D7").Value2 = 43465
$ws.Range("E7").Value2 = 43373
$ws.Range("D8").Value2 = 1291900
$ws.Range("E8").Value2 = 1231200
$ws.Range("D9").Value2 = 738900
$ws.Range("E9").Value2 = 693500
$ws.Range("D10").Value2 = 553000
$ws.Range("E10").Value2 = 537700
$ws.Range("D12").Value2 = "NA"
$ws.Range("E12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("E13").Value2 = 0
$ws.Range("D14").Value2 = 0
$ws.Range("E14").Value2 = 0
$ws.Range("D15").Value2 = 143700
$ws.Range("E15").Value2 = 81500
$ws.Range("D17").Value2 = 1225100
$ws.Range("E17").Value2 = 1162800
$ws.Range("D18").Value2 = 66800
$ws.Range("E18").Value2 = 68400
$ws.Range("D20").Value2 = 14500
$ws.Range("E20").Value2 = 5200
$ws.Range("D21").Value2 = 225100
$ws.Range("E21").Value2 = 155100
$ws.Range("D22").Value2 = 55500
$ws.Range("E22").Value2 = 49800
$ws.Range("D23").Value2 = 25800
$ws.Range("E23").Value2 = 23800
$ws.Range("D24").Value2 = 6800
$ws.Range("E24").Value2 = 9800
$ws.Range("D25").Value2 = 0
$ws.Range("E25").Value2 = 0
$ws.Range("D26").Value2 = 19100
$ws.Range("E26").Value2 = 14000
$ws.Range("D27").Value2 = 13100
$ws.Range("E27").Value2 = 4600
$ws.Range("D28").Value2 = 0
$ws.Range("E28").Value2 = 0
$ws.Range("D29").Value2 = 0
$ws.Range("E29").Value2 = 0
$ws.Range("D30").Value2 = 0
$ws.Range("E30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("E31").Value2 = 0
$ws.Range("D32").Value2 = -14500
$ws.Range("E32").Value2 = -5200
$ws.Range("D33").Value2 = 13100
$ws.Range("E33").Value2 = 4600
$ws.Range("D34").Value2 = 0
$ws.Range("E34").Value2 = 0
$ws.Range("D35").Value2 = 13100
$ws.Range("E35").Value2 = 4600
$ws.Range("D38").Value2 = 43465
$ws.Range("E38").Value2 = 43373
$ws.Range("D41").Value2 = 273100
$ws.Range("E41").Value2 = 403500
$ws.Range("D42").Value2 = "NA"
$ws.Range("E42").Value2 = 122100
$ws.Range("D43").Value2 = 526400
$ws.Range("E43").Value2 = 509600
$ws.Range("D44").Value2 = 90400
$ws.Range("E44").Value2 = 94600
$ws.Range("D45").Value2 = 226100
$ws.Range("E45").Value2 = 149600
$ws.Range("D46").Value2 = 1115900
$ws.Range("E46").Value2 = 1279400
$ws.Range("D47").Value2 = 35500
$ws.Range("E47").Value2 = 94600
$ws.Range("D48").Value2 = 5298500
$ws.Range("E48").Value2 = 5361900
$ws.Range("D49").Value2 = 513800
$ws.Range("E49").Value2 = 441500
$ws.Range("D50").Value2 = 0
$ws.Range("E50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("E51").Value2 = 0
$ws.Range("D52").Value2 = 140100
$ws.Range("E52").Value2 = 136800
$ws.Range("D53").Value2 = 0
$ws.Range("E53").Value2 = 0
$ws.Range("D54").Value2 = 7103800
$ws.Range("E54").Value2 = 7314100
$ws.Range("D57").Value2 = 506600
$ws.Range("E57").Value2 = 634400
$ws.Range("D58").Value2 = 626700
$ws.Range("E58").Value2 = 668800
$ws.Range("D59").Value2 = 1046200
$ws.Range("E59").Value2 = 921100
$ws.Range("D60").Value2 = 2179500
$ws.Range("E60").Value2 = 2224300
$ws.Range("D61").Value2 = 3380800
$ws.Range("E61").Value2 = 3464300
$ws.Range("D62").Value2 = 565800
$ws.Range("E62").Value2 = 576600
$ws.Range("D63").Value2 = 0
$ws.Range("E63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("E64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("E65").Value2 = 0
$ws.Range("D66").Value2 = 5948200
$ws.Range("E66").Value2 = 6098700
$ws.Range("D68").Value2 = 0
$ws.Range("E68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("E69").Value2 = 0
$ws.Range("D70").Value2 = 42000
$ws.Range("E70").Value2 = 42000
$ws.Range("D71").Value2 = 0
$ws.Range("E71").Value2 = 0
$ws.Range("D72").Value2 = 381300
$ws.Range("E72").Value2 = 328600
$ws.Range("D73").Value2 = 0
$ws.Range("E73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("E74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("E75").Value2 = 0
$ws.Range("D76").Value2 = 1113600
$ws.Range("E76").Value2 = 1173400
$ws.Range("D77").Value2 = 0
$ws.Range("E77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("E80").Value2 = 43373
$ws.Range("D81").Value2 = 13100
$ws.Range("E81").Value2 = 4600
$ws.Range("D83").Value2 = 143700
$ws.Range("E83").Value2 = 81500
$ws.Range("D84").Value2 = 0
$ws.Range("E84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("E85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("E86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("E87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("E88").Value2 = 0
$ws.Range("D89").Value2 = 240500
$ws.Range("E89").Value2 = 148200
$ws.Range("D91").Value2 = -241300
$ws.Range("E91").Value2 = -17500
$ws.Range("F91").Value2 = -119800
$ws.Range("G91").Value2 = -52100
$ws.Range("H91").Value2 = 1700
$ws.Range("I91").Value2 = -87300
$ws.Range("J91").Value2 = -88400
$ws.Range("D92").Value2 = 0
$ws.Range("E92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("E93").Value2 = 0
$ws.Range("D94").Value2 = -256800
$ws.Range("E94").Value2 = -43000
$ws.Range("D96").Value2 = -13500
$ws.Range("E96").Value2 = -16500
$ws.Range("D97").Value2 = 0
$ws.Range("E97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("E98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("E99").Value2 = 0
$ws.Range("D100").Value2 = -117000
$ws.Range("E100").Value2 = -91800
$ws.Range("D101").Value2 = 2800
$ws.Range("E101").Value2 = -7300
$ws.Range("D102").Value2 = -130400
$ws.Range("E102").Value2 = 6100
